$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 443.33334
$ws.Range("J121").Value = 441.15384
$ws.Range("L121").Value = 1323.46152
$ws.Range("N121").Value = -4817.46152

$ws.Range("H125").Value = 1090.2941
$ws.Range("I125").Value = 1790
$ws.Range("J125").Value = 875
$ws.Range("K125").Value = 16110
$ws.Range("L125").Value = 7875
$ws.Range("M125").Value = -13650
$ws.Range("N125").Value = -12795

$ws.Range("H129").Value = 1346.5333
$ws.Range("I129").Value = 642.9
$ws.Range("J129").Value = 1698.35
$ws.Range("K129").Value = 1928.7
$ws.Range("L129").Value = 5095.049999999999
$ws.Range("M129").Value = 3071.3
$ws.Range("N129").Value = -15095.05

$ws.Range("H138").Value = 4448.7627
$ws.Range("I138").Value = 2278.1562
$ws.Range("K138").Value = 6834.4686
$ws.Range("M138").Value = -1694.4686

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1401019.1
$ws.Range("I2").Value = 438.8889
$ws.Range("J2").Value = 2451454.2
$ws.Range("K2").Value = 438.8889
$ws.Range("L2").Value = 2451454.2
$ws.Range("M2").Value = -325.8889
$ws.Range("N2").Value = -2451680.2

$ws.Range("H32").Value = 12490.25
$ws.Range("I32").Value = 5468.596
$ws.Range("J32").Value = 58131
$ws.Range("K32").Value = 5468.596
$ws.Range("L32").Value = 58131
$ws.Range("M32").Value = -5181.596
$ws.Range("N32").Value = -58705

$ws.Range("H45").Value = 47888.316
$ws.Range("I45").Value = 85627.414
$ws.Range("J45").Value = 2601.4
$ws.Range("K45").Value = 85627.414
$ws.Range("L45").Value = 2601.4
$ws.Range("M45").Value = -85250.414
$ws.Range("N45").Value = -3355.4

$ws.Range("H116").Value = 1401019.1
$ws.Range("I116").Value = 438.8889
$ws.Range("J116").Value = 2451454.2
$ws.Range("K116").Value = 438.8889
$ws.Range("L116").Value = 2451454.2
$ws.Range("M116").Value = 1855.1111
$ws.Range("N116").Value = -2456042.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1401019.1
$ws.Range("I3").Value = 438.8889
$ws.Range("J3").Value = 2451454.2
$ws.Range("K3").Value = 438.8889
$ws.Range("L3").Value = 2451454.2
$ws.Range("M3").Value = -324.8889
$ws.Range("N3").Value = -2451682.2

$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H135").Value = 60000
$ws.Range("J135").Value = 60000
$ws.Range("L135").Value = 60000
$ws.Range("N135").Value = -70140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1567.2916
$ws.Range("I31").Value = 1193.6786
$ws.Range("J31").Value = 2090.35
$ws.Range("K31").Value = 1193.6786
$ws.Range("L31").Value = 2090.35
$ws.Range("M31").Value = -898.6786
$ws.Range("N31").Value = -2680.35

$ws.Range("H34").Value = 1567.2916
$ws.Range("I34").Value = 1193.6786
$ws.Range("J34").Value = 2090.35
$ws.Range("K34").Value = 1193.6786
$ws.Range("L34").Value = 2090.35
$ws.Range("M34").Value = -991.6786
$ws.Range("N34").Value = -2494.35

$ws.Range("H99").Value = 2149.75
$ws.Range("I99").Value = 1800
$ws.Range("J99").Value = 2499.5
$ws.Range("K99").Value = 1800
$ws.Range("L99").Value = 2499.5
$ws.Range("M99").Value = -302
$ws.Range("N99").Value = -5495.5

$ws.Range("H126").Value = 2149.75
$ws.Range("I126").Value = 1800
$ws.Range("J126").Value = 2499.5
$ws.Range("K126").Value = 5400
$ws.Range("L126").Value = 7498.5
$ws.Range("M126").Value = -2930
$ws.Range("N126").Value = -12438.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 10204591
$ws.Range("I5").Value = 524.75
$ws.Range("J5").Value = 17241878
$ws.Range("K5").Value = 1574.25
$ws.Range("L5").Value = 51725634
$ws.Range("M5").Value = -1462.25
$ws.Range("N5").Value = -51725858

$ws.Range("H68").Value = 984.6747
$ws.Range("I68").Value = 729.6905
$ws.Range("J68").Value = 1245.878
$ws.Range("K68").Value = 2189.0715
$ws.Range("L68").Value = 3737.634
$ws.Range("M68").Value = -1378.0715
$ws.Range("N68").Value = -5359.634

$ws.Range("H71").Value = 984.6747
$ws.Range("I71").Value = 729.6905
$ws.Range("J71").Value = 1245.878
$ws.Range("K71").Value = 6567.2145
$ws.Range("L71").Value = 11212.902
$ws.Range("M71").Value = -2511.2145
$ws.Range("N71").Value = -19324.902

$ws.Range("H131").Value = 1668687.9
$ws.Range("J131").Value = 2002136.2
$ws.Range("L131").Value = 6006408.6
$ws.Range("N131").Value = -6016488.6

$ws.Range("H135").Value = 10204591
$ws.Range("I135").Value = 524.75
$ws.Range("J135").Value = 17241878
$ws.Range("K135").Value = 4722.75
$ws.Range("L135").Value = 155176902
$ws.Range("M135").Value = -2187.75
$ws.Range("N135").Value = -155181972

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1182.4286
$ws.Range("I113").Value = 1499
$ws.Range("J113").Value = 945
$ws.Range("K113").Value = 1499
$ws.Range("L113").Value = 945
$ws.Range("M113").Value = 671
$ws.Range("N113").Value = -5285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1943.4445
$ws.Range("I7").Value = 1943.4445
$ws.Range("K7").Value = 1943.4445
$ws.Range("M7").Value = -1831.4445

$ws.Range("H40").Value = 3462.5
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 3616.6667
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 3616.6667
$ws.Range("M40").Value = -2864
$ws.Range("N40").Value = -3888.6667

$ws.Range("H126").Value = 1943.4445
$ws.Range("I126").Value = 1943.4445
$ws.Range("K126").Value = 5830.333500000001
$ws.Range("M126").Value = -3360.333500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 834.8333
$ws.Range("I122").Value = 699
$ws.Range("J122").Value = 1188
$ws.Range("K122").Value = 2097
$ws.Range("L122").Value = 3564
$ws.Range("M122").Value = 353
$ws.Range("N122").Value = -8464

$ws.Range("H126").Value = 2014.0667
$ws.Range("I126").Value = 2316.8333
$ws.Range("J126").Value = 803
$ws.Range("K126").Value = 6950.499899999999
$ws.Range("L126").Value = 2409
$ws.Range("M126").Value = -4480.499899999999
$ws.Range("N126").Value = -7349
